# Transactions.xlsx smoke-test update
#
# - Row 1 / column C gets the next amount in the series (stored as text,
#   same as the rest of the amount/card-number columns), replacing the old
#   "2109.02" value with "2634.02".
# - Row 1's trailing card-number cell (F1) is removed.
# - Row 2 (the second transaction) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# --- C1: replace the amount, keeping it text-typed like its neighbours ---
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "2634.02"
$ws.Range("C1").ClearFormats()

# --- F1: drop the trailing card-number cell on row 1 ---
$ws.Range("F1").Clear()

# --- Row 2: remove the whole second transaction row ---
$ws.Rows("2:2").Delete()
